# Add data for 2022-10-13: updates the "2022" (column I) totals across the
# Citywide, By Neighborhood and per-neighborhood sheets, plus a couple of
# incidental 2019 (column F) corrections picked up by the same refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("I2").Value = 5715
$ws.Range("I3").Value = 5967
$ws.Range("F4").Value = 1866
$ws.Range("I4").Value = 1360
$ws.Range("I5").Value = 547
$ws.Range("I6").Value = 6690
$ws.Range("F7").Value = 24055
$ws.Range("I7").Value = 20279

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("I2").Value = 66
$ws.Range("I3").Value = 52
$ws.Range("I7").Value = 220

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("I5").Value = 16
$ws.Range("I7").Value = 649

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("I2").Value = 183
$ws.Range("I3").Value = 292
$ws.Range("I7").Value = 795

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("I6").Value = 79
$ws.Range("I7").Value = 200

$ws = $wb.Worksheets.Item('New City')
$ws.Range("I3").Value = 144
$ws.Range("I4").Value = 20
$ws.Range("I7").Value = 473

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("I2").Value = 161
$ws.Range("I7").Value = 639
$ws.Range("I8").Value = 1220
$ws.Range("I9").Value = 99
$ws.Range("I11").Value = 300
$ws.Range("I15").Value = 229
$ws.Range("I17").Value = 32
$ws.Range("I18").Value = 145
$ws.Range("I19").Value = 559
$ws.Range("I20").Value = 492
$ws.Range("I22").Value = 52
$ws.Range("I24").Value = 56
$ws.Range("I25").Value = 104
$ws.Range("I29").Value = 1275
$ws.Range("I31").Value = 200
$ws.Range("I33").Value = 921
$ws.Range("I36").Value = 263
$ws.Range("I37").Value = 649
$ws.Range("I42").Value = 683
$ws.Range("I43").Value = 176
$ws.Range("I44").Value = 149
$ws.Range("I47").Value = 140
$ws.Range("I48").Value = 275
$ws.Range("I49").Value = 139
$ws.Range("I50").Value = 100
$ws.Range("I51").Value = 236
$ws.Range("F63").Value = 157
$ws.Range("I65").Value = 473
$ws.Range("I67").Value = 795
$ws.Range("I69").Value = 44
$ws.Range("I73").Value = 181
$ws.Range("I78").Value = 275
$ws.Range("I83").Value = 430
$ws.Range("I85").Value = 918
$ws.Range("I86").Value = 125
$ws.Range("I90").Value = 248
$ws.Range("I96").Value = 220
$ws.Range("I98").Value = 141
$ws.Range("F101").Value = 24055
$ws.Range("I101").Value = 20279

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("I2").Value = 148
$ws.Range("I7").Value = 430

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("I3").Value = 347
$ws.Range("I7").Value = 921

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("I6").Value = 85
$ws.Range("I7").Value = 139

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("I2").Value = 375
$ws.Range("I3").Value = 438
$ws.Range("I4").Value = 66
$ws.Range("I7").Value = 1275

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("I2").Value = 192
$ws.Range("I3").Value = 172
$ws.Range("I5").Value = 13
$ws.Range("I7").Value = 559

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("I3").Value = 44
$ws.Range("I7").Value = 149

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("I2").Value = 42
$ws.Range("I3").Value = 54
$ws.Range("I6").Value = 147
$ws.Range("I7").Value = 275

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("I2").Value = 251
$ws.Range("I3").Value = 354
$ws.Range("I4").Value = 47
$ws.Range("I5").Value = 32
$ws.Range("I7").Value = 918

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("I6").Value = 209
$ws.Range("I7").Value = 683

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("I2").Value = 65
$ws.Range("I7").Value = 275

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range("I2").Value = 20
$ws.Range("I7").Value = 56

$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Range("I6").Value = 14
$ws.Range("I7").Value = 44

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("I3").Value = 182
$ws.Range("I4").Value = 38

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("I6").Value = 167
$ws.Range("I7").Value = 492

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("I2").Value = 43
$ws.Range("I7").Value = 145

$ws = $wb.Worksheets.Item('Burnside')
$ws.Range("I6").Value = 4
$ws.Range("I7").Value = 32

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("I6").Value = 80
$ws.Range("I7").Value = 263

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("I2").Value = 37
$ws.Range("I4").Value = 7
$ws.Range("I7").Value = 104

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("I2").Value = 31
$ws.Range("I7").Value = 140

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("I3").Value = 52
$ws.Range("I7").Value = 229

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("I3").Value = 10
$ws.Range("I7").Value = 141

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("I3").Value = 21
$ws.Range("I7").Value = 100

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("I6").Value = 82
$ws.Range("I7").Value = 300

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("I6").Value = 30
$ws.Range("I7").Value = 99

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("I3").Value = 58
$ws.Range("I7").Value = 181

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("I2").Value = 54
$ws.Range("I7").Value = 161

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("I2").Value = 377
$ws.Range("I4").Value = 70
$ws.Range("I6").Value = 397
$ws.Range("I7").Value = 1220

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("I4").Value = 60
$ws.Range("I7").Value = 125

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("I4").Value = 13
$ws.Range("I7").Value = 248

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("I3").Value = 61
$ws.Range("I4").Value = 25
$ws.Range("I7").Value = 236

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("I2").Value = 36
$ws.Range("I7").Value = 176

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range("I2").Value = 22
$ws.Range("I7").Value = 52

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("I2").Value = 212
$ws.Range("I3").Value = 198
$ws.Range("I6").Value = 167
$ws.Range("I7").Value = 639
